$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: Set row 2 to reflect the (previously row 3 / FAPs->FAPs) data, with new TPM-derived values
$ws.Range("D2").Value = "FAPs"
$ws.Range("M2").Value = 0.1102156666666667
$ws.Range("N2").Value = 0.330647
$ws.Range("O2").Value = 0.8260949949157159
$ws.Range("P2").Value = 0.8260949949157158
$ws.Range("Q2").Value = 0.1787509277157778
$ws.Range("R2").Value = 1.608758349442
$ws.Range("S2").Value = 0.8260949949157159
$ws.Range("T2").Value = 0.8260949949157158

# Step 2: Set row 3 to reflect the (previously row 4 / FAPs->Resolving-Mac) data, with new TPM-derived values
$ws.Range("D3").Value = "Resolving-Mac"
$ws.Range("M3").Value = 0.023202
$ws.Range("N3").Value = 0.069606
$ws.Range("O3").Value = 0.1739050050842842
$ws.Range("P3").Value = 0.1739050050842842
$ws.Range("Q3").Value = 0.037629668724
$ws.Range("R3").Value = 0.338667018516
$ws.Range("S3").Value = 0.1739050050842842
$ws.Range("T3").Value = 0.1739050050842842

# Step 3: Remove old row 4 (FAPs->ECs pair no longer present), shifting everything up
$ws.Rows.Item(4).Delete()
